{"js": "const replacements = [\n  [\"305\u00d79=2745\", \"892\u00d73=2676\"],\n  [\"868\u00d79=7812\", \"228\u00d73=684\"],\n  [\"380\u00d77=2660\", \"338\u00d76=2028\"],\n  [\"559\u00d79=5031\", \"910\u00d74=3640\"],\n  [\"451\u00d75=2255\", \"591\u00d78=4728\"],\n  [\"793\u00d75=3965\", \"917\u00d79=8253\"],\n  [\"748\u00d72=1496\", \"222\u00d72=444\"],\n  [\"758\u00d74=3032\", \"259\u00d76=1554\"],\n  [\"984\u00d77=6888\", \"976\u00d72=1952\"],\n  [\"296\u00d73=888\", \"588\u00d72=1176\"],\n  [\"822\u00d75=4110\", \"977\u00d79=8793\"],\n  [\"220\u00d74=880\", \"866\u00d77=6062\"],\n  [\"425\u00d77=2975\", \"213\u00d73=639\"],\n  [\"969\u00d72=1938\", \"759\u00d79=6831\"],\n  [\"900\u00d78=7200\", \"425\u00d76=2550\"],\n  [\"118\u00d72=236\", \"816\u00d78=6528\"],\n  [\"315\u00d74=1260\", \"178\u00d79=1602\"],\n  [\"197\u00d72=394\", \"877\u00d77=6139\"],\n  [\"685\u00d73=2055\", \"843\u00d77=5901\"],\n  [\"566\u00d79=5094\", \"141\u00d77=987\"],\n  [\"670\u00d72=1340\", \"866\u00d79=7794\"],\n  [\"672\u00d74=2688\", \"972\u00d73=2916\"],\n  [\"163\u00d75=815\", \"188\u00d78=1504\"],\n  [\"751\u00d78=6008\", \"932\u00d75=4660\"],\n  [\"833\u00d75=4165\", \"695\u00d73=2085\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"305\u00d79=2745\", \"892\u00d73=2676\"),\n    @(\"868\u00d79=7812\", \"228\u00d73=684\"),\n    @(\"380\u00d77=2660\", \"338\u00d76=2028\"),\n    @(\"559\u00d79=5031\", \"910\u00d74=3640\"),\n    @(\"451\u00d75=2255\", \"591\u00d78=4728\"),\n    @(\"793\u00d75=3965\", \"917\u00d79=8253\"),\n    @(\"748\u00d72=1496\", \"222\u00d72=444\"),\n    @(\"758\u00d74=3032\", \"259\u00d76=1554\"),\n    @(\"984\u00d77=6888\", \"976\u00d72=1952\"),\n    @(\"296\u00d73=888\", \"588\u00d72=1176\"),\n    @(\"822\u00d75=4110\", \"977\u00d79=8793\"),\n    @(\"220\u00d74=880\", \"866\u00d77=6062\"),\n    @(\"425\u00d77=2975\", \"213\u00d73=639\"),\n    @(\"969\u00d72=1938\", \"759\u00d79=6831\"),\n    @(\"900\u00d78=7200\", \"425\u00d76=2550\"),\n    @(\"118\u00d72=236\", \"816\u00d78=6528\"),\n    @(\"315\u00d74=1260\", \"178\u00d79=1602\"),\n    @(\"197\u00d72=394\", \"877\u00d77=6139\"),\n    @(\"685\u00d73=2055\", \"843\u00d77=5901\"),\n    @(\"566\u00d79=5094\", \"141\u00d77=987\"),\n    @(\"670\u00d72=1340\", \"866\u00d79=7794\"),\n    @(\"672\u00d74=2688\", \"972\u00d73=2916\"),\n    @(\"163\u00d75=815\", \"188\u00d78=1504\"),\n    @(\"751\u00d78=6008\", \"932\u00d75=4660\"),\n    @(\"833\u00d75=4165\", \"695\u00d73=2085\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.Execute(\n        $oldText,    # FindText\n        $true,       # MatchCase\n        $false,      # MatchWholeWord\n        $false,      # MatchWildcards\n        $false,      # MatchSoundsLike\n        $false,      # MatchAllWordForms\n        $true,       # Forward\n        1,           # Wrap (wdFindContinue)\n        $false,      # Format\n        $newText,    # ReplaceWith\n        2            # Replace (wdReplaceAll)\n    )\n}\n"}
